$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 136, pushing old rows 136-166
# down to 139-169 (matches the dimension change A1:T166 -> A1:T169).
$ws.Rows("136:138").Insert()

# --- New row 136 ---
$ws.Range("A136").Value = 4
$ws.Range("B136").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C136").Value = "Los Lagos"
$ws.Range("D136").Value = 44551
$ws.Range("E136").Value = 10
$ws.Range("F136").Value = "Fruta"
$ws.Range("G136").Value = 100109
$ws.Range("H136").Value = "Uva"
$ws.Range("I136").Value = 100109001
$ws.Range("J136").Value = "Uva"
$ws.Range("K136").Value = "Flame Seedless"
$ws.Range("L136").Value = "Primera"
$ws.Range("M136").Value = 300
$ws.Range("N136").Value = 15000
$ws.Range("O136").Value = 16000
$ws.Range("P136").Value = 15500
$ws.Range("Q136").Value = "`$/bandeja 8 kilos"
$ws.Range("R136").Value = "Provincia de Copiapó"
$ws.Range("S136").Value = 1938
$ws.Range("T136").Value = 8

# --- New row 137 ---
$ws.Range("A137").Value = 4
$ws.Range("B137").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C137").Value = "Los Lagos"
$ws.Range("D137").Value = 44551
$ws.Range("E137").Value = 10
$ws.Range("F137").Value = "Fruta"
$ws.Range("G137").Value = 100109
$ws.Range("H137").Value = "Uva"
$ws.Range("I137").Value = 100109001
$ws.Range("J137").Value = "Uva"
$ws.Range("K137").Value = "Red Globe"
$ws.Range("L137").Value = "Primera"
$ws.Range("M137").Value = 300
$ws.Range("N137").Value = 19000
$ws.Range("O137").Value = 20000
$ws.Range("P137").Value = 19500
$ws.Range("Q137").Value = "`$/bandeja 8 kilos"
$ws.Range("R137").Value = "Provincia del Elquí"
$ws.Range("S137").Value = 2438
$ws.Range("T137").Value = 8

# --- New row 138 ---
$ws.Range("A138").Value = 4
$ws.Range("B138").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C138").Value = "Los Lagos"
$ws.Range("D138").Value = 44551
$ws.Range("E138").Value = 10
$ws.Range("F138").Value = "Fruta"
$ws.Range("G138").Value = 100109
$ws.Range("H138").Value = "Uva"
$ws.Range("I138").Value = 100109001
$ws.Range("J138").Value = "Uva"
$ws.Range("K138").Value = "Superior Seedless"
$ws.Range("L138").Value = "Primera"
$ws.Range("M138").Value = 300
$ws.Range("N138").Value = 17000
$ws.Range("O138").Value = 18000
$ws.Range("P138").Value = 17500
$ws.Range("Q138").Value = "`$/bandeja 8 kilos"
$ws.Range("R138").Value = "Provincia de Limarí"
$ws.Range("S138").Value = 2188
$ws.Range("T138").Value = 8
